# Apply edits to the case-locations-and-outbreaks public exposure sites sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in row 16, column D ("caught but" -> "caught bus")
$ws.Range("D16").Value = "Case caught bus from Melbourne Airport to Broadmeadows Railway Station"

# Delete rows from bottom to top so earlier row numbers stay valid.
# Row 38: West Melbourne / Kebab Kingz - removed entirely
$ws.Rows.Item(38).Delete()
# Row 22: South Melbourne / Stowe Australia - removed entirely
$ws.Rows.Item(22).Delete()
# Row 18: Melbourne / Exford Hotel - removed entirely
$ws.Rows.Item(18).Delete()
